$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 7.311265211180753

# Row 3
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "5c5882fc5bfe7600011197cb"
$ws.Range("E3").Value = "Colleen"
$ws.Range("G3").Value = 6.075952185643782
$ws.Range("H3").Value = "White"

# Row 4
$ws.Range("C4").Value = 19
$ws.Range("D4").Value = "60b45e9961dd412bfb6780f8"
$ws.Range("E4").Value = "Jewel"
$ws.Range("G4").Value = 6.068676626552405
$ws.Range("H4").Value = "Black or African American"

# Row 5
$ws.Range("G5").Value = 5.477047804629725

# Row 6
$ws.Range("G6").Value = 5.249471932023906

# Row 7
$ws.Range("G7").Value = 4.260356005502568

# Row 8
$ws.Range("G8").Value = 1.260598627945096

# Row 9
$ws.Range("G9").Value = 1.185192640848691

# Row 10
$ws.Range("C10").Value = 33
$ws.Range("D10").Value = "60cb36ee9f58331a33cf5506"
$ws.Range("E10").Value = "Shaniek"
$ws.Range("G10").Value = 0.3346982378612178

# Row 11
$ws.Range("C11").Value = 30
$ws.Range("D11").Value = "60d5775a99b502eec8cf56b4"
$ws.Range("E11").Value = "Shadaisia"
$ws.Range("G11").Value = 0.26099946291021
$ws.Range("H11").Value = "Black or African American"

# Row 12
$ws.Range("C12").Value = 32
$ws.Range("D12").Value = "6036f9b3b1842f8b659b18c7"
$ws.Range("E12").Value = "Kellie"
$ws.Range("G12").Value = 0.07698541627100014
$ws.Range("H12").Value = "White"

# Row 13
$ws.Range("G13").Value = 0.04919117767745862

# Row 14
$ws.Range("G14").Value = 13.19982871425305

# Row 15
$ws.Range("G15").Value = 8.081433205567341

# Row 16
$ws.Range("C16").Value = 30
$ws.Range("D16").Value = "60c2341fe95d71ee52c043f0"
$ws.Range("E16").Value = "Matthew"
$ws.Range("G16").Value = 7.048241805590385

# Row 17
$ws.Range("C17").Value = 27
$ws.Range("D17").Value = "5ff8ad350d084e10f500e48a"
$ws.Range("E17").Value = "Drew"
$ws.Range("G17").Value = 7.000947600168775

# Row 18
$ws.Range("C18").Value = 26
$ws.Range("D18").Value = "5dd671942b033b5ec8bc97b4"
$ws.Range("E18").Value = "Juan"
$ws.Range("G18").Value = 5.239313832273305
$ws.Range("H18").Value = "Hispanic"

# Row 19
$ws.Range("C19").Value = 22
$ws.Range("D19").Value = "60db4fde6193c50664c9c478"
$ws.Range("E19").Value = "Edosagbe"
$ws.Range("G19").Value = 5.186302527479196
$ws.Range("H19").Value = "Black or African American"

# Row 20
$ws.Range("G20").Value = 5.018452747422359

# Row 21
$ws.Range("G21").Value = 4.155549573790759

# Row 22
$ws.Range("G22").Value = 4.034175108618071

# Row 23
$ws.Range("G23").Value = 3.262000135003892

# Row 24
$ws.Range("G24").Value = 2.30063950806506

# Row 25
$ws.Range("G25").Value = 2.244526961475056
